$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H52").Value = 732.25
$ws.Range("I52").Value = 264.5
$ws.Range("J52").Value = 1200
$ws.Range("K52").Value = 793.5
$ws.Range("L52").Value = 3600
$ws.Range("M52").Value = -633.5
$ws.Range("N52").Value = -3920
$ws.Range("H111").Value = 1579
$ws.Range("I111").Value = 1294.4286
$ws.Range("J111").Value = 2575
$ws.Range("K111").Value = 3883.2858
$ws.Range("L111").Value = 7725
$ws.Range("M111").Value = -816.2857999999997
$ws.Range("N111").Value = -13859
$ws.Range("H137").Value = 4411.5625
$ws.Range("I137").Value = 1758.6
$ws.Range("K137").Value = 5275.799999999999
$ws.Range("M137").Value = -2725.799999999999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1435.8948
$ws.Range("I61").Value = 1440.6666
$ws.Range("K61").Value = 1440.6666
$ws.Range("M61").Value = -1228.6666
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H132").Value = 2240.8845
$ws.Range("I132").Value = 1930.56
$ws.Range("J132").Value = 9999
$ws.Range("K132").Value = 5791.68
$ws.Range("L132").Value = 29997
$ws.Range("M132").Value = -3261.68
$ws.Range("N132").Value = -35057
$ws.Range("H136").Value = 1435.8948
$ws.Range("I136").Value = 1440.6666
$ws.Range("K136").Value = 4321.9998
$ws.Range("M136").Value = -1771.9998
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4064.9722
$ws.Range("I105").Value = 3502.05
$ws.Range("K105").Value = 3502.05
$ws.Range("M105").Value = -1755.05
$ws.Range("H134").Value = 2352.8235
$ws.Range("I134").Value = 835.1667
$ws.Range("J134").Value = 5995.2
$ws.Range("K134").Value = 2505.5001
$ws.Range("L134").Value = 17985.6
$ws.Range("M134").Value = 29.4998999999998
$ws.Range("N134").Value = -23055.6
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 15003.5
$ws.Range("J4").Value = 30002
$ws.Range("L4").Value = 30002
$ws.Range("N4").Value = -30226
$ws.Range("H16").Value = 5000
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 5000
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 5000
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -5574
$ws.Range("H58").Value = 6840.923
$ws.Range("I58").Value = 5231
$ws.Range("K58").Value = 5231
$ws.Range("M58").Value = -5028
$ws.Range("H113").Value = 5000
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 5000
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -9340
$ws.Range("H132").Value = 3501.4783
$ws.Range("J132").Value = 4249.3335
$ws.Range("L132").Value = 12748.0005
$ws.Range("N132").Value = -17808.0005
$ws.Range("H136").Value = 6840.923
$ws.Range("I136").Value = 5231
$ws.Range("K136").Value = 15693
$ws.Range("M136").Value = -13143
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4685985.5
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H69").Value = 3674.3333
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 3674.3333
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("I132").Value = 2000
$ws.Range("J132").Value = 10749.5
$ws.Range("K132").Value = 18000
$ws.Range("L132").Value = 96745.5
$ws.Range("M132").Value = -15470
$ws.Range("N132").Value = -101805.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 20002
$ws.Range("I5").Value = 20002
$ws.Range("K5").Value = 20002
$ws.Range("M5").Value = -19890
$ws.Range("H70").Value = 7000
$ws.Range("I70").Value = 7000
$ws.Range("K70").Value = 7000
$ws.Range("M70").Value = -6730
$ws.Range("H73").Value = 7000
$ws.Range("I73").Value = 7000
$ws.Range("K73").Value = 7000
$ws.Range("M73").Value = -6064
$ws.Range("H134").Value = 69000
$ws.Range("J134").Value = 69000
$ws.Range("L134").Value = 207000
$ws.Range("N134").Value = -212070
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 5057500
$ws.Range("I2").Value = 20000000
$ws.Range("J2").Value = 76666.336
$ws.Range("K2").Value = 20000000
$ws.Range("L2").Value = 76666.336
$ws.Range("M2").Value = -19999888
$ws.Range("N2").Value = -76890.336
$ws.Range("H46").Value = 3056.9167
$ws.Range("I46").Value = 1891.2667
$ws.Range("K46").Value = 1891.2667
$ws.Range("M46").Value = -1703.2667
$ws.Range("H61").Value = 7138.4287
$ws.Range("I61").Value = 8992.333000000001
$ws.Range("K61").Value = 8992.333000000001
$ws.Range("M61").Value = -8790.333000000001
$ws.Range("H113").Value = 7138.4287
$ws.Range("I113").Value = 8992.333000000001
$ws.Range("K113").Value = 8992.333000000001
$ws.Range("M113").Value = -6822.333000000001
$ws.Range("H122").Value = 5000
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 15000
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -19900
$ws.Range("H127").Value = 82374.75
$ws.Range("J127").Value = 82374.75
$ws.Range("L127").Value = 82374.75
$ws.Range("N127").Value = -92294.75
$ws.Range("H132").Value = 5669.7036
$ws.Range("I132").Value = 3783.5715
$ws.Range("K132").Value = 11350.7145
$ws.Range("M132").Value = -8820.7145
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 3333666.8
$ws.Range("J2").Value = 3333666.8
$ws.Range("L2").Value = 3333666.8
$ws.Range("N2").Value = -3333890.8
$ws.Range("H81").Value = 2229.0588
$ws.Range("I81").Value = 2229.0588
$ws.Range("K81").Value = 4458.1176
$ws.Range("M81").Value = -3397.1176
$ws.Range("H84").Value = 2229.0588
$ws.Range("I84").Value = 2229.0588
$ws.Range("K84").Value = 22290.588
$ws.Range("M84").Value = -16986.588
$ws.Range("H107").Value = 985.8
$ws.Range("I107").Value = 798.1429000000001
$ws.Range("J107").Value = 1150
$ws.Range("K107").Value = 2394.4287
$ws.Range("L107").Value = 3450
$ws.Range("M107").Value = -474.4287000000004
$ws.Range("N107").Value = -7290
$ws.Range("H113").Value = 1800
$ws.Range("J113").Value = 1974.25
$ws.Range("L113").Value = 5922.75
$ws.Range("N113").Value = -10262.75
$ws.Range("H132").Value = 1704.6897
$ws.Range("I132").Value = 1215.4546
$ws.Range("K132").Value = 3646.3638
$ws.Range("M132").Value = -1116.3638
$ws.Range("H136").Value = 55596.58
$ws.Range("I136").Value = 2784.7646
$ws.Range("K136").Value = 8354.293799999999
$ws.Range("M136").Value = -5804.293799999999
